$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source rule table value in C10 (row for "R20") changes from 18 to 1.
$ws.Range("C10").Value = 1
